# Update "想去人数" (F column) counts across sheets, matching the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 318
$ws.Range("F4").Value = 2982
$ws.Range("F7").Value = 2325
$ws.Range("F8").Value = 1688
$ws.Range("F10").Value = 856
$ws.Range("F13").Value = 2666
$ws.Range("F15").Value = 1531
$ws.Range("F16").Value = 7085
$ws.Range("F18").Value = 7234
$ws.Range("F21").Value = 5489
$ws.Range("F22").Value = 3115
$ws.Range("F24").Value = 236
$ws.Range("F25").Value = 187
$ws.Range("F26").Value = 1900
$ws.Range("F28").Value = 303
$ws.Range("F29").Value = 878
$ws.Range("F31").Value = 285
$ws.Range("F32").Value = 40
$ws.Range("F33").Value = 2428
$ws.Range("F34").Value = 1203
$ws.Range("F35").Value = 2735
$ws.Range("F38").Value = 170
$ws.Range("F39").Value = 394
$ws.Range("F40").Value = 1086
$ws.Range("F43").Value = 527

# --- Sheet: 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 157
$ws.Range("F14").Value = 49

# --- Sheet: 本地生活 (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 58

# --- Sheet: 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 318
$ws.Range("F5").Value = 2982
$ws.Range("F7").Value = 2325
$ws.Range("F8").Value = 1688
$ws.Range("F11").Value = 856
$ws.Range("F14").Value = 2666
$ws.Range("F15").Value = 1531
$ws.Range("F19").Value = 7085
$ws.Range("F21").Value = 7234
$ws.Range("F23").Value = 5489
$ws.Range("F24").Value = 3115
$ws.Range("F27").Value = 236
$ws.Range("F28").Value = 49
$ws.Range("F29").Value = 1900
$ws.Range("F32").Value = 303
$ws.Range("F33").Value = 878
$ws.Range("F35").Value = 285
$ws.Range("F36").Value = 40
$ws.Range("F37").Value = 2428
$ws.Range("F38").Value = 1203
$ws.Range("F40").Value = 2735
$ws.Range("F43").Value = 170
$ws.Range("F45").Value = 394
$ws.Range("F46").Value = 1086
$ws.Range("F49").Value = 527
